# Adds "pair_words"-column (D) carrier duplicates for the practice rows,
# and adds "unique_video"/"unique_audio" pair_kind/carrier entries for the
# generic word rows, matching the upstream edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Practice rows (2-5): add column D with the carrier word (mirrors column K)
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Generic rows 6-9: add column J pair_kind value
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# Rows 14-21: add pair_kind (C) and carrier (D) values
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "can"
$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "can"
$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "do"
$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "do"
$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "look"
$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "look"
$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "where"
$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "where"
